# Add an "Address" column: insert a new column at F, pushing the existing
# "District" column (F) one position right to G, then fill the new F column
# (rows 2 header + data rows 4-38) with the school address extracted from
# column B's second line.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at F. This shifts the old F column ("District"
# header in F2/"NAME" in F3 and the district values in F4:F38) to column G.
$ws.Columns.Item(6).Insert()

# New header for the inserted column.
$ws.Range("F2").Value = "Address"

# New "Address" values for data rows 4-38 (row 3's F cell stays blank).
$addresses = @(
    @(4,  "GunariSmt G P H Govt. High School Kavaloor"),
    @(5,  "S M V H S TavarageraKustagi"),
    @(6,  "S G High School"),
    @(7,  "G H S KesarahattiGangavati"),
    @(8,  "Govt. J R College"),
    @(9,  "Govt. P U CollegeYelburga"),
    @(10, "Girls P U CollegeKustagi"),
    @(11, "BenakallamathGovt. Girls High School Talakal"),
    @(12, "S J A High School Gondabal"),
    @(13, "G G H S Hanamasagar"),
    @(14, "G H S Navali Gangavathi"),
    @(15, "G G H S Yelburga"),
    @(16, "G H S LingadahalliKushtagi"),
    @(17, "Govt. High School DanapurGangavathi"),
    @(18, "S K C P U C CollegeGangavathi"),
    @(19, "G H S NavalahalliKuishtagi"),
    @(20, "G H S H V KuntaYelburga"),
    @(21, "G P U C Irakalgada"),
    @(22, "G H S NilogalKushtagi"),
    @(23, "S R S M H S Hitnal"),
    @(24, "G H S Horatatnal"),
    @(25, "G H S KuknoorYelburga"),
    @(26, "Viveka Bharathi High SchoolGangavati"),
    @(27, "Govt. High SchoolVajrabandiYelaburga"),
    @(28, "G H S MukkumpiGangavatti"),
    @(29, "Boys G P U CollegeKustagi"),
    @(30, "G H S HulihyderGangavathi"),
    @(31, "G H S YarageraKustagi"),
    @(32, "G H S IslampurGangavathi"),
    @(33, "G H S HirejanthkalGangavathi"),
    @(34, "G H S MataladinniYelburga"),
    @(35, "G H S Chilakamukhi"),
    @(36, "G J CollegeTalakalYelaburga"),
    @(37, "Govt High SchoolKolur"),
    @(38, "V N C P U C Munirabad")
)

foreach ($pair in $addresses) {
    $row = $pair[0]
    $value = $pair[1]
    $ws.Cells.Item($row, 6).Value = $value
}
